# Gros update refonte code
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LIST")
$ws2 = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------------
# Sheet "LIST" (sheet1)
# ---------------------------------------------------------------------------

# Row 3: drop D3, add (empty, text-formatted) B3
$ws1.Range("D3").Clear() | Out-Null
$ws1.Range("B3").NumberFormat = "@"

# Row 4: drop D4
$ws1.Range("D4").Clear() | Out-Null

# Row 5: value moves from AD.SEC.014.FON.01 to AD.DEP.001.FON.01, and gets the
# grey "s=2" text style; drop D5
$ws1.Range("A5").Value = "AD.DEP.001.FON.01"
$ws1.Range("A5").NumberFormat = "@"
$ws1.Range("D5").Clear() | Out-Null

# New rows 6-11 (styled like the rest of column A)
$ws1.Range("A6").Value = "RO.ACT.001"
$ws1.Range("A6").NumberFormat = "@"
$ws1.Range("A7").Value = "RO.ACT.003"
$ws1.Range("A7").NumberFormat = "@"
$ws1.Range("A8").Value = "RO.ACT.005"
$ws1.Range("A8").NumberFormat = "@"
$ws1.Range("A9").Value = "RO.FOU.001"
$ws1.Range("A9").NumberFormat = "@"
$ws1.Range("A10").Value = "MP.CPT.001"
$ws1.Range("A10").NumberFormat = "@"
$ws1.Range("A11").Value = "RT.ART.001"
$ws1.Range("A11").NumberFormat = "@"

# New rows 12-13 (no special style)
$ws1.Range("A12").Value = "RT.MAT.001"
$ws1.Range("A13").Value = "RO.ORG.001"

# Row 14: the value that used to live in A5
$ws1.Range("A14").Value = "AD.SEC.014.FON.01"

# Column width bookkeeping: drop the old (now unused) E/F/G/L/M overrides and
# give H/I the width that L/M used to carry
$ws1.Columns.Item(8).ColumnWidth = 20.83
$ws1.Columns.Item(9).ColumnWidth = 20.83

# ---------------------------------------------------------------------------
# Sheet "Feuil1" (sheet2)
# ---------------------------------------------------------------------------

$ws2.Range("F22").Value = "RO.ACT.001"
$ws2.Range("F22").NumberFormat = "@"
$ws2.Range("F23").Value = "RO.ACT.003"
$ws2.Range("F23").NumberFormat = "@"
$ws2.Range("F24").Value = "RO.ACT.005"
$ws2.Range("F24").NumberFormat = "@"
$ws2.Range("F25").Value = "RO.FOU.001"
$ws2.Range("F25").NumberFormat = "@"
$ws2.Range("F26").Value = "MP.CPT.001"
$ws2.Range("F26").NumberFormat = "@"
$ws2.Range("F27").Value = "RT.ART.001"
$ws2.Range("F27").NumberFormat = "@"
$ws2.Range("F28").Value = "RT.MAT.001"
$ws2.Range("F29").Value = "RO.ORG.001"

# ---------------------------------------------------------------------------
# Selections / active window state (applied last so sheet "LIST" ends up the
# active tab, matching the saved file)
# ---------------------------------------------------------------------------

$ws2.Activate()
$ws2.Range("B11:B14").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10

$ws1.Activate()
$ws1.Range("F1:F1048576").Select() | Out-Null
